# lite-Locale.xlsx: add a new data row (row 52, 2019-04-14 counts) and move
# the active cell to B52 (region/selection bug fix).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data (date serial 43569 == 2019-04-14)
$rowValues = @(43569, 6, 5, 1, 381, 15, 34, 3, 1, 3, 1, 5, 4, 2, 1, 13, 2, 0, 0)

for ($col = 1; $col -le $rowValues.Length; $col++) {
    $ws.Cells.Item(52, $col).Value = $rowValues[$col - 1]
}

# Match the date formatting already used by column A (reuse the existing
# style instead of inventing a new number format)
$ws.Range("A51").Copy()
$ws.Range("A52").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Region / selection bug fix: move the active cell onto the new row, at B52
$ws.Range("B52").Select()
